# Applies the "Penalty Reward System" week-shift edits described in the
# commit diff:
#   - Sheet "Forecast Comparison": shift every Week_Start_Date (col B,
#     rows 2-17) forward by one week.
#   - Sheet "Summary": update the dependent summary statistics (historical
#     range end date, total/forecast counts, max/min forecast week) to
#     reflect the shifted data.
#
# All of the target cells hold plain text (e.g. "2025-01-05", "35",
# "136 units") rather than real Excel dates/numbers, so every value is
# written with a leading apostrophe to force text entry and avoid Excel's
# automatic date/number parsing from changing the cell type.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Forecast Comparison ---------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekStartDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

foreach ($row in $weekStartDates.Keys) {
    $ws1.Range("B$row").Value = "'" + $weekStartDates[$row]
}

# ---- Sheet: Summary -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").Value  = "'2023-10-15 to 2025-01-05"
$ws2.Range("B8").Value  = "'138 units"
$ws2.Range("B9").Value  = "'32"
$ws2.Range("B10").Value = "'16"
$ws2.Range("B11").Value = "'8"
$ws2.Range("B13").Value = "'2025-01-12"
$ws2.Range("B15").Value = "'2025-02-02"
